# Mise à jour des backlogs et assignation des sprints
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Assigner les user stories aux sprints (colonne G : "Assignée au sprint")
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 2
$ws.Range("G5").Value = 3

# Déplacer la sélection active sur C5
$ws.Range("C5").Select()
